# Inserts one new daily price record for "Haba" (Región Metropolitana,
# 2021-10-26) into the weekly/daily consolidated sheet. The new record is
# inserted as row 58, pushing the existing rows 58..160 down to 59..161
# (matching the way Excel's Rows.Insert shifts a block down and keeps the
# row's formatting, e.g. the date-number-format on column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 58 - everything below shifts
# down by one (old row 58 becomes row 59, ..., old row 160 becomes row 161).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new record's data.
$ws.Range("A58").Value = 9
$ws.Range("B58").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C58").Value = "Metropolitana"
$ws.Range("D58").Value = 44495
$ws.Range("E58").Value = 13
$ws.Range("F58").Value = 100112026
$ws.Range("G58").Value = "Haba"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 170
$ws.Range("K58").Value = 7000
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = 7471
$ws.Range("N58").Value = "$/saco 25 kilos"
$ws.Range("O58").Value = "Región Metropolitana"
$ws.Range("P58").Value = 299
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
